$wb = $excel.ActiveWorkbook

$oldId = "7e802949-0df4-4339-9423-8e93f0f0289a"
$newId = "8b10ff1d-7258-479b-9cee-88ff74bfa152"

$oldMd = "$oldId.md"
$newMd = "$newId.md"

$oldZhXlf = "$oldId.166dcceb795b225a4ef056c04c3f77a7b1a66fc3.zh-cn.xlf"
$newZhXlf = "$newId.7929c32eb86e154314d1e8ead1699068515dae77.zh-cn.xlf"

$oldDeXlf = "$oldId.166dcceb795b225a4ef056c04c3f77a7b1a66fc3.de-de.xlf"
$newDeXlf = "$newId.7929c32eb86e154314d1e8ead1699068515dae77.de-de.xlf"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("G2").Value = "2016-08-16 18:52:14"

# B2 carries a hyperlink whose external target stays the same (stale in the
# source commit); only the cell text / hyperlink display text changes.
$oldOverviewAddress = $wsOverview.Hyperlinks.Item(1).Address
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Range("B2").Hyperlinks.Add($wsOverview.Range("B2"), $oldOverviewAddress, [Type]::Missing, [Type]::Missing, "e2e\$newMd")

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$oldZhAddress = $wsZh.Hyperlinks.Item(1).Address
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("A2").Hyperlinks.Add($wsZh.Range("A2"), $oldZhAddress, [Type]::Missing, [Type]::Missing, $newMd)

$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = "2016-08-16 18:52:09"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$oldDeAddress = $wsDe.Hyperlinks.Item(1).Address
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("A2").Hyperlinks.Add($wsDe.Range("A2"), $oldDeAddress, [Type]::Missing, [Type]::Missing, $newMd)

$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = "2016-08-16 18:52:14"
